# Update LR-pairs data with new TPM-derived values
# (Sema6a-Plxna2 ligand-receptor pair recomputed after TPM update)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 43.841169
$ws.Range("H2").Value = 131.523507
$ws.Range("I2").Value = 0.6105408572336042
$ws.Range("J2").Value = 0.6105408572336042
$ws.Range("M2").Value = 30.61919899999999
$ws.Range("N2").Value = 91.85759699999998
$ws.Range("O2").Value = 0.6951390881735714
$ws.Range("P2").Value = 0.6951390881735714
$ws.Range("Q2").Value = 1342.381478003631
$ws.Range("R2").Value = 12081.43330203268
$ws.Range("S2").Value = 0.4244108147900783
$ws.Range("T2").Value = 0.4244108147900783

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 43.841169
$ws.Range("H3").Value = 131.523507
$ws.Range("I3").Value = 0.6105408572336042
$ws.Range("J3").Value = 0.6105408572336042
$ws.Range("O3").Value = 0.09263417906992544
$ws.Range("P3").Value = 0.09263417906992545
$ws.Range("Q3").Value = 178.885647964728
$ws.Range("R3").Value = 1609.970831682552
$ws.Range("S3").Value = 0.05655695109848347
$ws.Range("T3").Value = 0.05655695109848348

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 43.841169
$ws.Range("H4").Value = 131.523507
$ws.Range("I4").Value = 0.6105408572336042
$ws.Range("J4").Value = 0.6105408572336042
$ws.Range("M4").Value = 9.348075333333332
$ws.Range("N4").Value = 28.04422599999999
$ws.Range("O4").Value = 0.2122267327565031
$ws.Range("P4").Value = 0.2122267327565031
$ws.Range("Q4").Value = 409.830550513398
$ws.Range("R4").Value = 3688.474954620581
$ws.Range("S4").Value = 0.1295730913450424
$ws.Range("T4").Value = 0.1295730913450424

$ws.Range("I5").Value = 0.0635739353967235
$ws.Range("J5").Value = 0.06357393539672351
$ws.Range("M5").Value = 30.61919899999999
$ws.Range("N5").Value = 91.85759699999998
$ws.Range("O5").Value = 0.6951390881735714
$ws.Range("P5").Value = 0.6951390881735714
$ws.Range("Q5").Value = 139.77848058694
$ws.Range("R5").Value = 1258.00632528246
$ws.Range("S5").Value = 0.04419272748328391
$ws.Range("T5").Value = 0.04419272748328392

$ws.Range("I6").Value = 0.0635739353967235
$ws.Range("J6").Value = 0.06357393539672351
$ws.Range("O6").Value = 0.09263417906992544
$ws.Range("P6").Value = 0.09263417906992545
$ws.Range("S6").Value = 0.005889119315719956
$ws.Range("T6").Value = 0.005889119315719958

$ws.Range("I7").Value = 0.0635739353967235
$ws.Range("J7").Value = 0.06357393539672351
$ws.Range("M7").Value = 9.348075333333332
$ws.Range("N7").Value = 28.04422599999999
$ws.Range("O7").Value = 0.2122267327565031
$ws.Range("P7").Value = 0.2122267327565031
$ws.Range("Q7").Value = 42.67452478118666
$ws.Range("R7").Value = 384.07072303068
$ws.Range("S7").Value = 0.01349208859771963
$ws.Range("T7").Value = 0.01349208859771964

$ws.Range("G8").Value = 23.400872
$ws.Range("H8").Value = 70.202616
$ws.Range("I8").Value = 0.3258852073696723
$ws.Range("J8").Value = 0.3258852073696723
$ws.Range("M8").Value = 30.61919899999999
$ws.Range("N8").Value = 91.85759699999998
$ws.Range("O8").Value = 0.6951390881735714
$ws.Range("P8").Value = 0.6951390881735714
$ws.Range("Q8").Value = 716.515956541528
$ws.Range("R8").Value = 6448.643608873752
$ws.Range("S8").Value = 0.2265355459002092
$ws.Range("T8").Value = 0.2265355459002092

$ws.Range("G9").Value = 23.400872
$ws.Range("H9").Value = 70.202616
$ws.Range("I9").Value = 0.3258852073696723
$ws.Range("J9").Value = 0.3258852073696723
$ws.Range("O9").Value = 0.09263417906992544
$ws.Range("P9").Value = 0.09263417906992545
$ws.Range("Q9").Value = 95.48285883206401
$ws.Range("R9").Value = 859.3457294885761
$ws.Range("S9").Value = 0.03018810865572201
$ws.Range("T9").Value = 0.03018810865572201

$ws.Range("G10").Value = 23.400872
$ws.Range("H10").Value = 70.202616
$ws.Range("I10").Value = 0.3258852073696723
$ws.Range("J10").Value = 0.3258852073696723
$ws.Range("M10").Value = 9.348075333333332
$ws.Range("N10").Value = 28.04422599999999
$ws.Range("O10").Value = 0.2122267327565031
$ws.Range("P10").Value = 0.2122267327565031
$ws.Range("Q10").Value = 218.7531143216907
$ws.Range("R10").Value = 1968.778028895216
$ws.Range("S10").Value = 0.06916155281374105
$ws.Range("T10").Value = 0.06916155281374105
